$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new Price (column D) values are plain numeric-looking strings
# (e.g. "281.74"). Excel auto-converts such strings to numbers on assignment,
# but the original workbook stores every Price cell as text (inlineStr). To
# preserve that, force these specific cells to Text format before assigning,
# then restore their style to "Normal" so the logical number format (General)
# is unchanged afterwards. Each cell is handled individually (not as a single
# multi-area Union range) because Range operations on Union ranges in this
# runtime only reliably affect the first area.
$textForceCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D17", "D19", "D20", "D24", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '21.097.17'
$ws.Range("E2").Value = '  +3.40%  '
$ws.Range("D3").Value = '1.535.95'
$ws.Range("E3").Value = '  +5.10%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").Value = '0.9663'
$ws.Range("E5").Value = '  +1.66%  '
$ws.Range("D6").Value = '281.74'
$ws.Range("E6").Value = '  +2.59%  '
$ws.Range("D7").Value = '0.3619'
$ws.Range("E7").Value = '  -0.93%  '
$ws.Range("D8").Value = '0.3170'
$ws.Range("E8").Value = '  +3.56%  '
$ws.Range("D9").Value = '40.67'
$ws.Range("E9").Value = '  +2.23%  '
$ws.Range("D10").Value = '1.093'
$ws.Range("D11").Value = '0.06808'
$ws.Range("E11").Value = '  +3.52%  '
$ws.Range("D12").Value = '1.007'
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("D13").Value = '5.654'
$ws.Range("E13").Value = '  +4.43%  '
$ws.Range("D14").Value = '18.70'
$ws.Range("E14").Value = '  +4.02%  '
$ws.Range("D15").Value = '6.345'
$ws.Range("E15").Value = '  +3.39%  '
$ws.Range("E16").Value = '  +1.97%  '
$ws.Range("D17").Value = '0.9663'
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").Value = '1.528.17'
$ws.Range("E18").Value = '  +4.49%  '
$ws.Range("D19").Value = '0.06080'
$ws.Range("E19").Value = '  +4.45%  '
$ws.Range("D20").Value = '72.23'
$ws.Range("E20").Value = '  +4.57%  '
$ws.Range("E21").Value = '  +4.98%  '
$ws.Range("E22").Value = '  +4.01%  '
$ws.Range("E23").Value = '  +4.26%  '
$ws.Range("D24").Value = '2.317'
$ws.Range("E24").Value = '  +3.36%  '
$ws.Range("D25").Value = '21.165.40'
$ws.Range("E25").Value = '  +3.65%  '
$ws.Range("D26").Value = '148.20'
$ws.Range("E26").Value = '  +4.59%  '
$ws.Range("D27").Value = '2.214'
$ws.Range("E27").Value = '  +6.92%  '
$ws.Range("D28").Value = '17.63'
$ws.Range("E28").Value = '  +3.14%  '
$ws.Range("D29").Value = '1.694.94'
$ws.Range("E29").Value = '  +4.93%  '
$ws.Range("D30").Value = '118.45'
$ws.Range("E30").Value = '  +5.10%  '
$ws.Range("D31").Value = '4.020'
$ws.Range("E31").Value = '  +4.41%  '
$ws.Range("D32").Value = '0.8506'
$ws.Range("E32").Value = '  +7.78%  '
$ws.Range("D33").Value = '5.166'
$ws.Range("E33").Value = '  +5.76%  '
$ws.Range("D34").Value = '0.07995'
$ws.Range("E34").Value = '  +1.48%  '
$ws.Range("D35").Value = '1.508'
$ws.Range("E35").Value = '  -0.54%  '
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").Value = '1.211'
$ws.Range("E36").Value = '  +5.96%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '4.950'
$ws.Range("E37").Value = '  +5.95%  '
$ws.Range("D38").Value = '0.05855'
$ws.Range("E38").Value = '  +2.59%  '
$ws.Range("D39").Value = '0.02097'
$ws.Range("E39").Value = '  +3.56%  '
$ws.Range("D40").Value = '10.63'
$ws.Range("E40").Value = '  +3.12%  '
$ws.Range("D41").Value = '7.722'
$ws.Range("E41").Value = '  +3.48%  '
$ws.Range("D42").Value = '0.9666'
$ws.Range("E42").Value = '  +0.84%  '
$ws.Range("D43").Value = '0.1909'
$ws.Range("E43").Value = '  +2.87%  '
$ws.Range("D44").Value = '0.5427'
$ws.Range("E44").Value = '  +3.33%  '
$ws.Range("D45").Value = '12.49'
$ws.Range("E45").Value = '  +5.28%  '
$ws.Range("D46").Value = '3.578'
$ws.Range("E46").Value = '  +2.62%  '
$ws.Range("D47").Value = '0.5432'
$ws.Range("E47").Value = '  +5.54%  '
$ws.Range("D48").Value = '121.34'
$ws.Range("E48").Value = '  +3.59%  '
$ws.Range("D49").Value = '1.866'
$ws.Range("E49").Value = '  +6.92%  '
$ws.Range("D50").Value = '0.06576'
$ws.Range("E50").Value = '  +2.47%  '
$ws.Range("D51").Value = '0.9933'
$ws.Range("E51").Value = '  +0.07%  '

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
